$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Séance TPA2-A4-D et TPA3 du 30/01/2024 : deux nouvelles lignes de journal
# de bord (lignes 7 et 8), jusque-là vides.
# ---------------------------------------------------------------------------

# -- Row 7 : groupe A2-4 (MPAL / TP) --
$ws.Range("A7").Value = 45321
$ws.Range("A7").NumberFormat = "m/d/yyyy"
$ws.Range("B7").Value = "MPAL"
$ws.Range("C7").Value = "TP"
$ws.Range("E7").Value = "x"

# -- Row 8 : groupe A3 (MPAL / TP) --
$ws.Range("A8").Value = 45321
$ws.Range("A8").NumberFormat = "m/d/yyyy"
$ws.Range("B8").Value = "MPAL"
$ws.Range("C8").Value = "TP"
$ws.Range("F8").Value = "x"

# -- New shared-string cells, entered in authoring (column-major) order so
#    they land in the sharedStrings table in the same order as the source --
$ws.Range("G7").Value = "Simple stack : correction #2.1 + indications de correction données pour finir jusqu'à #2.3.`nFinir #2.2 et #2.3 + préparer projet Questionsscore pour la prochaine fois."
$ws.Range("G8").Value = "Simple stack : correction #2.2 et #2.3.`nDebut questionsscore : #0 et debut kata sur #1. J'ai montré comment faire un Generate de Test."
$ws.Range("H8").Value = "Certains font un clone via Gitbash dans le même répertoir sans utiliser l'interface Get From VCS de IntelliJ et ça pose des problèmes… clarifier qu'il faut toujours passer par IntelliJ pour cloner le dépôt."
$ws.Range("I7").Value = "Certains oublient le testPeekOnEmptyStack. La démarche générale semble avoir été comprise.`nIl a fallu du temps à certains pour réactiver la licence sur le poste, avec la manipulation du proxy…"
$ws.Range("I8").Value = "Pas de problème particulier. Le groupe a bien compris.`nCertains testent assertNotThrown dans un cas nominal de peek ou pop, mais ce n'est pas utile il me semble car si une exception est levée, le test ne passe pas.`nIls ont des souvenirs de l'utilisation d'Eclipse dans lequel il faudrait faire ça pour que les tests passents sans se bloquer. `nOr dans IntelliJ les tests semblent bien tous indépendants."

# -- Formatting: vertically centre the whole of both new rows, and wrap the
#    long free-text cells so the multi-line comments are fully visible --
$ws.Range("A7:I7").VerticalAlignment = -4108
$ws.Range("G7:I7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 58

$ws.Range("A8:F8").VerticalAlignment = -4108
$ws.Range("G8").WrapText = $true
$ws.Range("I8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 87

# -- Leave the selection on the last cell touched, like the source edit --
$ws.Range("I8").Select() | Out-Null
